$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-307) holds a "Förändrad" (changed) date serial number.
# All of these cells currently hold 45182 and must be updated to 45184.
$ws.Range("C2:C307").Value = 45184
